# Highlight three specific list-item sentences in yellow.
# Both the paragraph mark (w:pPr/w:rPr) and the run text (w:r/w:rPr)
# need the <w:highlight w:val="yellow"/> applied, which corresponds to
# selecting the whole paragraph (including its end-of-paragraph mark)
# and setting Range.HighlightColorIndex = wdYellow (7).

$d = $word.ActiveDocument

$targets = @(
    "Para cada tipo de audio, canciones y podcast, informar el acumulado total de reproducciones en toda la plataforma.",
    "Informar el género de canción más escuchado en toda la plataforma y su número de reproducciones.",
    "Informar la categoría de podcast más escuchada en toda la plataforma y su número de reproducciones."
)

foreach ($t in $targets) {
    $rng = $d.Content
    $found = $rng.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        # Extend the range to include the paragraph mark so that the
        # paragraph's rPr (w:pPr/w:rPr) also gets the highlight, matching
        # how Word highlights an entire selected paragraph.
        $para = $rng.Paragraphs(1)
        $paraRange = $para.Range
        $paraRange.HighlightColorIndex = 7
    }
}
